$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns A (Id), Q (Ost), R (Nord) for rows 2-17.
# All other columns in each row are unchanged; only these three columns'
# values are being redistributed among the rows.
$updates = @(
    @{ Row = 2;  A = 111336414; Q = 625154.518985283;  R = 7209790.550712772 }
    @{ Row = 3;  A = 111336398; Q = 625080.5973707421; R = 7209863.719748351 }
    @{ Row = 4;  A = 111336417; Q = 625114.8910281583; R = 7209835.074288641 }
    @{ Row = 5;  A = 111336408; Q = 625192.9989858982; R = 7209753.182408583 }
    @{ Row = 6;  A = 111336411; Q = 625220.1037653659; R = 7209754.321236268 }
    @{ Row = 7;  A = 111336403; Q = 625077.0981605061; R = 7209815.690764531 }
    @{ Row = 8;  A = 111336404; Q = 625133.9709027896; R = 7209754.516890368 }
    @{ Row = 9;  A = 111336407; Q = 625184.0874587877; R = 7209753.231787121 }
    @{ Row = 10; A = 111336412; Q = 625184.8249035137; R = 7209765.975211025 }
    @{ Row = 11; A = 111336416; Q = 625105.0111505401; R = 7209827.879692691 }
    @{ Row = 12; A = 111336413; Q = 625162.5685057295; R = 7209780.71893465  }
    @{ Row = 13; A = 111336419; Q = 625180.6727454782; R = 7209784.021863313 }
    @{ Row = 14; A = 111336410; Q = 625259.9138955096; R = 7209755.994336623 }
    @{ Row = 15; A = 111336415; Q = 625148.1990682605; R = 7209799.60762905  }
    @{ Row = 16; A = 111336409; Q = 625220.7730415409; R = 7209758.586850428 }
    @{ Row = 17; A = 111336405; Q = 625163.9543035047; R = 7209747.724855823 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.A   # Column A
    $ws.Cells.Item($r, 17).Value = $u.Q  # Column Q
    $ws.Cells.Item($r, 18).Value = $u.R  # Column R
}
